$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1454.2858
$ws.Range("I12").Value = 196.66667
$ws.Range("J12").Value = 9000
$ws.Range("K12").Value = 196.66667
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = -26.66667000000001
$ws.Range("N12").Value = -9340
$ws.Range("H28").Value = 1628.2222
$ws.Range("I28").Value = 1093.4286
$ws.Range("K28").Value = 1093.4286
$ws.Range("M28").Value = -608.4286
$ws.Range("H51").Value = 5198.6665
$ws.Range("I51").Value = 4998
$ws.Range("K51").Value = 4998
$ws.Range("M51").Value = -4514
$ws.Range("H58").Value = 228.18182
$ws.Range("I58").Value = 228.18182
$ws.Range("K58").Value = 684.5454599999999
$ws.Range("M58").Value = -534.5454599999999
$ws.Range("H62").Value = 1999.5
$ws.Range("I62").Value = 1999.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1999.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1375.5
$ws.Range("N62").Value = $null
$ws.Range("H65").Value = 1999.5
$ws.Range("I65").Value = 1999.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9997.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6877.5
$ws.Range("N65").Value = $null
$ws.Range("H80").Value = 3855.3103
$ws.Range("I80").Value = 2100.1667
$ws.Range("J80").Value = 4313.174
$ws.Range("K80").Value = 6300.500100000001
$ws.Range("L80").Value = 12939.522
$ws.Range("M80").Value = -5302.500100000001
$ws.Range("N80").Value = -14935.522
$ws.Range("H83").Value = 3855.3103
$ws.Range("I83").Value = 2100.1667
$ws.Range("J83").Value = 4313.174
$ws.Range("K83").Value = 18901.5003
$ws.Range("L83").Value = 38818.566
$ws.Range("M83").Value = -13909.5003
$ws.Range("N83").Value = -48802.566
$ws.Range("H87").Value = 19904.762
$ws.Range("J87").Value = 19904.762
$ws.Range("L87").Value = 19904.762
$ws.Range("N87").Value = -22400.762
$ws.Range("H90").Value = 19904.762
$ws.Range("J90").Value = 19904.762
$ws.Range("L90").Value = 59714.28599999999
$ws.Range("N90").Value = -72194.28599999999
$ws.Range("H106").Value = 1674
$ws.Range("I106").Value = 896.6667
$ws.Range("K106").Value = 896.6667
$ws.Range("M106").Value = -265.6667
$ws.Range("H107").Value = 2238.2666
$ws.Range("I107").Value = 2508
$ws.Range("J107").Value = 1698.8
$ws.Range("K107").Value = 2508
$ws.Range("L107").Value = 1698.8
$ws.Range("M107").Value = -588
$ws.Range("N107").Value = -5538.8
$ws.Range("H132").Value = 15731.815
$ws.Range("I132").Value = 1373.3572
$ws.Range("J132").Value = 55935.5
$ws.Range("K132").Value = 4120.071599999999
$ws.Range("L132").Value = 167806.5
$ws.Range("M132").Value = -1590.071599999999
$ws.Range("N132").Value = -172866.5
$ws.Range("H135").Value = 26324928
$ws.Range("I135").Value = 41672120
$ws.Range("K135").Value = 375049080
$ws.Range("M135").Value = -375046545
$ws.Range("H138").Value = 3747.56
$ws.Range("I138").Value = 2248.75
$ws.Range("K138").Value = 6746.25
$ws.Range("M138").Value = -1606.25
$ws.Range("H141").Value = 7266.5483
$ws.Range("I141").Value = 4856.269
$ws.Range("K141").Value = 14568.807
$ws.Range("M141").Value = -9388.807000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15879525
$ws.Range("I32").Value = 17864172
$ws.Range("K32").Value = 17864172
$ws.Range("M32").Value = -17863885
$ws.Range("H61").Value = 3459.5
$ws.Range("I61").Value = 2895.76
$ws.Range("K61").Value = 2895.76
$ws.Range("M61").Value = -2683.76
$ws.Range("H74").Value = 2224.2727
$ws.Range("I74").Value = 2501.8
$ws.Range("J74").Value = 1993
$ws.Range("K74").Value = 2501.8
$ws.Range("L74").Value = 1993
$ws.Range("M74").Value = -1627.8
$ws.Range("N74").Value = -3741
$ws.Range("H77").Value = 2224.2727
$ws.Range("I77").Value = 2501.8
$ws.Range("J77").Value = 1993
$ws.Range("K77").Value = 12509
$ws.Range("L77").Value = 9965
$ws.Range("M77").Value = -8141
$ws.Range("N77").Value = -18701
$ws.Range("H88").Value = 1452.2354
$ws.Range("I88").Value = 1406.8572
$ws.Range("J88").Value = 1484
$ws.Range("K88").Value = 1406.8572
$ws.Range("L88").Value = 1484
$ws.Range("M88").Value = -1000.8572
$ws.Range("N88").Value = -2296
$ws.Range("H91").Value = 1452.2354
$ws.Range("I91").Value = 1406.8572
$ws.Range("J91").Value = 1484
$ws.Range("K91").Value = 1406.8572
$ws.Range("L91").Value = 1484
$ws.Range("M91").Value = -2.857199999999921
$ws.Range("N91").Value = -4292
$ws.Range("H102").Value = 3431.0454
$ws.Range("I102").Value = 3431.0454
$ws.Range("K102").Value = 3431.0454
$ws.Range("M102").Value = -1809.0454
$ws.Range("H110").Value = 2276
$ws.Range("I110").Value = 2562.842
$ws.Range("K110").Value = 2562.842
$ws.Range("M110").Value = -517.8420000000001
$ws.Range("H132").Value = 8327.439
$ws.Range("I132").Value = 8147.423
$ws.Range("J132").Value = 8639.467000000001
$ws.Range("K132").Value = 24442.269
$ws.Range("L132").Value = 25918.401
$ws.Range("M132").Value = -21912.269
$ws.Range("N132").Value = -30978.401
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
$ws.Range("H136").Value = 3459.5
$ws.Range("I136").Value = 2895.76
$ws.Range("K136").Value = 8687.280000000001
$ws.Range("M136").Value = -6137.280000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 1500
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -1275
$ws.Range("N64").Value = -4450
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 1500
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -720
$ws.Range("N67").Value = -5560
$ws.Range("H86").Value = 85148.086
$ws.Range("I86").Value = 168100
$ws.Range("J86").Value = 2196.1667
$ws.Range("K86").Value = 168100
$ws.Range("L86").Value = 2196.1667
$ws.Range("M86").Value = -166977
$ws.Range("N86").Value = -4442.1667
$ws.Range("H89").Value = 85148.086
$ws.Range("I89").Value = 168100
$ws.Range("J89").Value = 2196.1667
$ws.Range("K89").Value = 840500
$ws.Range("L89").Value = 10980.8335
$ws.Range("M89").Value = -834884
$ws.Range("N89").Value = -22212.8335
$ws.Range("H99").Value = 20010.908
$ws.Range("I99").Value = 21807.25
$ws.Range("J99").Value = 2047.5
$ws.Range("K99").Value = 21807.25
$ws.Range("L99").Value = 2047.5
$ws.Range("M99").Value = -20309.25
$ws.Range("N99").Value = -5043.5
$ws.Range("H105").Value = 2374.8333
$ws.Range("I105").Value = 1649.8
$ws.Range("K105").Value = 1649.8
$ws.Range("M105").Value = 97.20000000000005
$ws.Range("H107").Value = 9302.200000000001
$ws.Range("I107").Value = 10503.667
$ws.Range("J107").Value = 7500
$ws.Range("K107").Value = 10503.667
$ws.Range("L107").Value = 7500
$ws.Range("M107").Value = -8583.666999999999
$ws.Range("N107").Value = -11340
$ws.Range("H133").Value = 99796.414
$ws.Range("I133").Value = 97777
$ws.Range("J133").Value = 99980
$ws.Range("K133").Value = 97777
$ws.Range("L133").Value = 99980
$ws.Range("M133").Value = -92717
$ws.Range("N133").Value = -110100
$ws.Range("H134").Value = 4191.5557
$ws.Range("I134").Value = 4302.7334
$ws.Range("J134").Value = 3635.6667
$ws.Range("K134").Value = 12908.2002
$ws.Range("L134").Value = 10907.0001
$ws.Range("M134").Value = -10373.2002
$ws.Range("N134").Value = -15977.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4195.5615
$ws.Range("I31").Value = 5119
$ws.Range("J31").Value = 3169.5186
$ws.Range("K31").Value = 5119
$ws.Range("L31").Value = 3169.5186
$ws.Range("M31").Value = -4824
$ws.Range("N31").Value = -3759.5186
$ws.Range("H34").Value = 4195.5615
$ws.Range("I34").Value = 5119
$ws.Range("J34").Value = 3169.5186
$ws.Range("K34").Value = 5119
$ws.Range("L34").Value = 3169.5186
$ws.Range("M34").Value = -4917
$ws.Range("N34").Value = -3573.5186
$ws.Range("H62").Value = 14500
$ws.Range("I62").Value = 14000
$ws.Range("K62").Value = 14000
$ws.Range("M62").Value = -13376
$ws.Range("H65").Value = 14500
$ws.Range("I65").Value = 14000
$ws.Range("K65").Value = 70000
$ws.Range("M65").Value = -66880
$ws.Range("H133").Value = 77431.25
$ws.Range("J133").Value = 78242
$ws.Range("L133").Value = 78242
$ws.Range("N133").Value = -83302
$ws.Range("H134").Value = 2313.5652
$ws.Range("I134").Value = 2632.457
$ws.Range("K134").Value = 7897.370999999999
$ws.Range("M134").Value = -5362.370999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 40.57143
$ws.Range("J2").Value = 22.5
$ws.Range("L2").Value = 135
$ws.Range("N2").Value = -361
$ws.Range("H7").Value = 16549.77
$ws.Range("I7").Value = 16549.77
$ws.Range("K7").Value = 49649.31
$ws.Range("M7").Value = -49537.31
$ws.Range("H8").Value = 901
$ws.Range("I8").Value = 901
$ws.Range("K8").Value = 2703
$ws.Range("M8").Value = -2564
$ws.Range("H38").Value = 801.57574
$ws.Range("I38").Value = 97.72727
$ws.Range("K38").Value = 293.18181
$ws.Range("M38").Value = 53.81818999999996
$ws.Range("H68").Value = 1515.9231
$ws.Range("I68").Value = 1200
$ws.Range("J68").Value = 2021.4
$ws.Range("K68").Value = 3600
$ws.Range("L68").Value = 6064.200000000001
$ws.Range("M68").Value = -2789
$ws.Range("N68").Value = -7686.200000000001
$ws.Range("H71").Value = 1515.9231
$ws.Range("I71").Value = 1200
$ws.Range("J71").Value = 2021.4
$ws.Range("K71").Value = 10800
$ws.Range("L71").Value = 18192.6
$ws.Range("M71").Value = -6744
$ws.Range("N71").Value = -26304.6
$ws.Range("H81").Value = 1470.75
$ws.Range("I81").Value = 1628.3334
$ws.Range("K81").Value = 4885.0002
$ws.Range("M81").Value = -3762.0002
$ws.Range("H84").Value = 1470.75
$ws.Range("I84").Value = 1628.3334
$ws.Range("K84").Value = 14655.0006
$ws.Range("M84").Value = -9039.000599999999
$ws.Range("H117").Value = 1179.1111
$ws.Range("J117").Value = 1867
$ws.Range("L117").Value = 5601
$ws.Range("N117").Value = -12485
$ws.Range("H121").Value = 1454.3334
$ws.Range("I121").Value = 917.2727
$ws.Range("J121").Value = 1722.8636
$ws.Range("K121").Value = 2751.8181
$ws.Range("L121").Value = 5168.5908
$ws.Range("M121").Value = -1441.8181
$ws.Range("N121").Value = -7788.5908
$ws.Range("H122").Value = 4241.609
$ws.Range("I122").Value = 881.6667
$ws.Range("J122").Value = 4745.6
$ws.Range("K122").Value = 7935.0003
$ws.Range("L122").Value = 42710.4
$ws.Range("M122").Value = -5485.0003
$ws.Range("N122").Value = -47610.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 244.6875
$ws.Range("I2").Value = 140.8
$ws.Range("J2").Value = 417.83334
$ws.Range("K2").Value = 140.8
$ws.Range("L2").Value = 417.83334
$ws.Range("M2").Value = -27.80000000000001
$ws.Range("N2").Value = -643.83334
$ws.Range("H33").Value = 7545.778
$ws.Range("J33").Value = 8364
$ws.Range("L33").Value = 8364
$ws.Range("N33").Value = -8868
$ws.Range("H40").Value = 3633
$ws.Range("I40").Value = 2949.75
$ws.Range("K40").Value = 2949.75
$ws.Range("M40").Value = -2798.75
$ws.Range("H42").Value = 48831.668
$ws.Range("J42").Value = 48831.668
$ws.Range("L42").Value = 48831.668
$ws.Range("N42").Value = -49801.668
$ws.Range("H44").Value = 12767.333
$ws.Range("J44").Value = 14599.75
$ws.Range("L44").Value = 14599.75
$ws.Range("N44").Value = -15791.75
$ws.Range("H47").Value = 25930.5
$ws.Range("J47").Value = 25930.5
$ws.Range("L47").Value = 25930.5
$ws.Range("N47").Value = -27066.5
$ws.Range("H55").Value = 5998
$ws.Range("I55").Value = 5998
$ws.Range("K55").Value = 5998
$ws.Range("M55").Value = -5671
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H97").Value = 34500
$ws.Range("I97").Value = 500
$ws.Range("K97").Value = 500
$ws.Range("M97").Value = -4
$ws.Range("H115").Value = 48831.668
$ws.Range("J115").Value = 48831.668
$ws.Range("L115").Value = 48831.668
$ws.Range("N115").Value = -51181.668
$ws.Range("H132").Value = 9106.308000000001
$ws.Range("I132").Value = 8383.679
$ws.Range("K132").Value = 25151.037
$ws.Range("M132").Value = -22621.037

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7862.6665
$ws.Range("I16").Value = 2595.5
$ws.Range("K16").Value = 2595.5
$ws.Range("M16").Value = -2425.5
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H24").Value = 5000
$ws.Range("I24").Value = 5000
$ws.Range("K24").Value = 5000
$ws.Range("M24").Value = -4657
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H40").Value = 4682.846
$ws.Range("I40").Value = 4590
$ws.Range("J40").Value = 4992.3335
$ws.Range("K40").Value = 4590
$ws.Range("L40").Value = 4992.3335
$ws.Range("M40").Value = -4454
$ws.Range("N40").Value = -5264.3335
$ws.Range("H55").Value = 937.125
$ws.Range("I55").Value = 642.4286
$ws.Range("K55").Value = 642.4286
$ws.Range("M55").Value = -469.4286
$ws.Range("H94").Value = 62500
$ws.Range("I94").Value = 40000
$ws.Range("K94").Value = 40000
$ws.Range("M94").Value = -39324
$ws.Range("H100").Value = 195537.45
$ws.Range("I100").Value = 303701.16
$ws.Range("J100").Value = 6251
$ws.Range("K100").Value = 303701.16
$ws.Range("L100").Value = 6251
$ws.Range("M100").Value = -303160.16
$ws.Range("N100").Value = -7333
$ws.Range("H130").Value = 69000
$ws.Range("J130").Value = 65000
$ws.Range("L130").Value = 65000
$ws.Range("N130").Value = -75040
$ws.Range("H136").Value = 2137.2
$ws.Range("I136").Value = 1529.0588
$ws.Range("K136").Value = 4587.1764
$ws.Range("M136").Value = -2037.1764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 63022
$ws.Range("I96").Value = 106799.6
$ws.Range("J96").Value = 8300
$ws.Range("K96").Value = 106799.6
$ws.Range("L96").Value = 8300
$ws.Range("M96").Value = -105426.6
$ws.Range("N96").Value = -11046
$ws.Range("H113").Value = 1343.5454
$ws.Range("J113").Value = 1618.5
$ws.Range("L113").Value = 4855.5
$ws.Range("N113").Value = -9195.5
$ws.Range("H126").Value = 2224.75
$ws.Range("I126").Value = 2224.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6674.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4204.25
$ws.Range("N126").Value = $null
$ws.Range("H128").Value = 66000
$ws.Range("J128").Value = 66000
$ws.Range("L128").Value = 66000
$ws.Range("N128").Value = -75960
$ws.Range("H129").Value = 74993
$ws.Range("J129").Value = 74993
$ws.Range("L129").Value = 74993
$ws.Range("N129").Value = -84993
$ws.Range("H132").Value = 3142.25
$ws.Range("I132").Value = 3370.75
$ws.Range("K132").Value = 10112.25
$ws.Range("M132").Value = -7582.25
$ws.Range("H133").Value = 108135
$ws.Range("J133").Value = 108135
$ws.Range("L133").Value = 108135
$ws.Range("N133").Value = -118255
$ws.Range("H136").Value = 3434.16
$ws.Range("I136").Value = 2857.1875
$ws.Range("J136").Value = 4459.8887
$ws.Range("K136").Value = 8571.5625
$ws.Range("L136").Value = 13379.6661
$ws.Range("M136").Value = -6021.5625
$ws.Range("N136").Value = -18479.6661
